# Regenerate save_data to use K (strikeouts) instead of Strike# in column G.
# Column G (header "K") values for rows 2-8 are recalculated/rewritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2 = 4
    3 = 1
    4 = 3
    5 = 3
    6 = 1
    7 = 4
    8 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
